$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2020-06-24"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = 5

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "10"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "100"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "Yes"
